$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2200.2292173303581
$ws.Range("B1").Value = 1413.0883371171401
$ws.Range("C1").Value = 1433.2414832486022
$ws.Range("A2").Value = 2223.1192368732227
$ws.Range("B2").Value = 1464.6353980171634
$ws.Range("C2").Value = 1430.9469458087237
$ws.Range("A3").Value = 2339.9137069531034
$ws.Range("B3").Value = 1581.8914897206105
$ws.Range("C3").Value = 1450.1273560749194
$ws.Range("A4").Value = 2316.96250606476
$ws.Range("B4").Value = 1770.3572927336343
$ws.Range("C4").Value = 1699.6232104240059
$ws.Range("A5").Value = 2418.6539934030511
$ws.Range("B5").Value = 1664.5656950833943
$ws.Range("C5").Value = 1621.8655132393455
$ws.Range("A6").Value = 2356.7084402705214
$ws.Range("B6").Value = 1770.6284616285568
$ws.Range("C6").Value = 1780.4372524670989
$ws.Range("A7").Value = 1988.9848474206719
$ws.Range("B7").Value = 1563.0572054573668
$ws.Range("C7").Value = 1480.331633107181
$ws.Range("A8").Value = 2131.6712526944602
$ws.Range("B8").Value = 1649.3121108864455
$ws.Range("C8").Value = 1628.0931546045269
$ws.Range("A9").Value = 2467.015372537443
$ws.Range("B9").Value = 1784.3089756160193
$ws.Range("C9").Value = 1509.8937671913866
$ws.Range("A10").Value = 2107.594404518145
$ws.Range("B10").Value = 1366.7866019700048
$ws.Range("C10").Value = 1283.1297772827584
$ws.Range("A11").Value = 1966.0038839175647
$ws.Range("B11").Value = 1412.2680734948408
$ws.Range("C11").Value = 1294.508315442717
$ws.Range("A12").Value = 2783.9537336872691
$ws.Range("B12").Value = 2266.437411002204
$ws.Range("C12").Value = 2032.7945288481224
$ws.Range("A13").Value = 2311.1853410557019
$ws.Range("B13").Value = 1738.1827811731157
$ws.Range("C13").Value = 1814.3571781233311
$ws.Range("A14").Value = 2589.0481324808579
$ws.Range("B14").Value = 1918.3549811213027
$ws.Range("C14").Value = 1698.6671106768149
$ws.Range("A15").Value = 2504.3188623771234
$ws.Range("B15").Value = 2043.7271445454799
$ws.Range("C15").Value = 1841.8802398524401
$ws.Range("A16").Value = 2201.5211166450968
$ws.Range("B16").Value = 1506.1457634936228
$ws.Range("C16").Value = 1266.566756430175
$ws.Range("A17").Value = 2221.5448154951432
$ws.Range("B17").Value = 1685.6426137539731
$ws.Range("C17").Value = 1559.9267688263026
$ws.Range("A18").Value = 2483.9072682735809
$ws.Range("B18").Value = 2057.2767182058769
$ws.Range("C18").Value = 1913.4119084784868
$ws.Range("A19").Value = 1737.5582140970034
$ws.Range("B19").Value = 1921.0538878036466
$ws.Range("C19").Value = 1884.8258167361525
$ws.Range("A20").Value = 2347.276345455467
$ws.Range("B20").Value = 1821.6045095213988
$ws.Range("C20").Value = 1764.6830779821748
$ws.Range("A21").Value = 2579.110070587506
$ws.Range("B21").Value = 1903.0651052615717
$ws.Range("C21").Value = 1809.8882081724635
$ws.Range("A22").Value = 2444.3691932947804
$ws.Range("B22").Value = 1886.2358905041021
$ws.Range("C22").Value = 1641.6541911979325
